# Update the SOLAR expense (row 22) on Sheet1 from 332.05 to 222.
# Dependent formulas (E28 = SUM(E10:E27) and G39 = E7-E28-E36) will
# recalculate automatically to reflect the new total.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("E22").Value = 222

# Reflect the author's final cursor/selection position on the sheet.
$ws.Range("K32").Select()
